# "Added Test Data for UK Market"
#
# Duplicate the "Poland" sheet (it's the template for a per-market test
# data sheet) to a new sheet placed right after it, rename the copy to
# "UK", and swap in the UK-specific market label / Jira reference. Also
# activate + select the new sheet the way the author would have left it
# selected after finishing the edit (Excel automatically clears
# tabSelected on the previously-active sheet).

$wb = $excel.ActiveWorkbook

$poland = $wb.Worksheets.Item("Poland")

# Copy "Poland" to a new sheet immediately after it.
$poland.Copy($null, $poland)
$uk = $wb.Worksheets.Item($wb.Worksheets.Count)
$uk.Name = "UK"

# Update the market-specific cells (order matters for shared-string index
# allocation: the Jira/ticket reference is entered before the market name).
$uk.Range("B4").Value = "NGC-2741/T3343/T3342"
$uk.Range("B2").Value = "UK Market"

# Leave the new "UK" sheet active/selected, cursor on B4.
$uk.Activate()
$uk.Range("B4").Select()
